# Delete comments & Complete this report.
#
# The review comment left on slide 1 ("提出する人は記入忘れずに！ / それは私
# かもしれない． / Write the submission date.") has been addressed, so remove
# it. Deleting the last remaining comment on the deck removes
# ppt/comments/comment1.xml (and its Content_Types override / slide
# relationship) while leaving ppt/commentAuthors.xml untouched.

$p = $ppt.ActivePresentation

foreach ($s in $p.Slides) {
    for ($i = $s.Comments.Count; $i -ge 1; $i--) {
        $s.Comments.Item($i).Delete()
    }
}
